# Updates the cryptos list data (columns B-E) to the latest scraped values.
# Column D (Price) values are forced to text via a leading apostrophe so Excel
# does not reinterpret dotted/decimal strings as numbers (avoiding float drift
# and loss of formats like "76.316.83" or trailing zeros like "3.80").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.299.02"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "'2.968.41"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'637.05"
$ws.Range("E5").Value = "  +7.06%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'199.35"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "'2.967.64"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "'0.432"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "'3.514.68"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").Value = "'28.82"
$ws.Range("E15").Value = "  +5.90%  "
$ws.Range("D16").Value = "'76.267.49"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "'0.0000187"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "'2.965.09"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "'13.35"
$ws.Range("E19").Value = "  +6.44%  "
$ws.Range("D20").Value = "'8.71"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'370.61"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").Value = "'72.63"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").Value = "'3.123.55"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").Value = "'9.58"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "'8.23"
$ws.Range("E31").Value = "  +7.23%  "
$ws.Range("D32").Value = "'513.60"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'1.96"
$ws.Range("E34").Value = "  +9.15%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'20.18"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").Value = "'163.27"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +12.52%  "
$ws.Range("D39").Value = "'19.96"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "  +15.91%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").Value = "'181.71"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E44").Value = "  +6.91%  "
$ws.Range("D45").Value = "'4.89"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "'0.701"
$ws.Range("E48").Value = "  +7.92%  "
$ws.Range("D49").Value = "'0.581"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("D51").Value = "'3.80"
$ws.Range("E51").Value = "  +2.61%  "

Write-Host "Applied cryptos list update"
